# Add two new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing formatting used by the other header/data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto the new
# header cells so they match the bold/centered/bordered look.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data values (rows 2-64) ---
$data = @(
    @(2, 6, 8),
    @(3, 9, 9),
    @(4, 9, 10),
    @(5, 8, 9),
    @(6, 6, 7),
    @(7, 7, 7),
    @(8, 9, 9),
    @(9, 8, 8),
    @(10, 8, 9),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 9, 9),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 6, 6),
    @(17, 7, 7),
    @(18, 6, 6),
    @(19, 6, 6),
    @(20, 5, 5),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 7, 8),
    @(24, 8, 8),
    @(25, 7, 7),
    @(26, 8, 8),
    @(27, 6, 6),
    @(28, 8, 8),
    @(29, 9, 9),
    @(30, 8, 9),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 8, 8),
    @(34, 11, 11),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 6, 6),
    @(38, 7, 7),
    @(39, 7, 7),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 9, 9),
    @(43, 6, 7),
    @(44, 6, 6),
    @(45, 8, 8),
    @(46, 6, 6),
    @(47, 7, 7),
    @(48, 7, 7),
    @(49, 8, 8),
    @(50, 7, 7),
    @(51, 8, 8),
    @(52, 7, 7),
    @(53, 10, 10),
    @(54, 8, 8),
    @(55, 9, 9),
    @(56, 8, 8),
    @(57, 9, 9),
    @(58, 8, 8),
    @(59, 8, 8),
    @(60, 8, 8),
    @(61, 5, 5),
    @(62, 6, 6),
    @(63, 6, 6),
    @(64, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
